# Reposition a handful of connector / rectangle shapes across slides 2, 3, 4 and 9.
# Offsets in the source OOXML are expressed in EMUs (914400 EMU = 1 inch,
# 12700 EMU = 1 point). The PowerPoint COM object model exposes shape
# positions in points via .Left / .Top, so EMU values are divided by 12700.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# The COM host stores Shape.Left/.Top as single-precision (f32) points and
# truncates (floors) when converting back to EMU on write. Nudging the
# point value up by half an EMU keeps the round-trip (pts -> f32 -> EMU)
# landing on the intended integer EMU instead of one unit short.
$HalfEmuInPoints = 0.5 / 12700.0

function Move-Shape($slideIndex, $shapeId, $xEmu, $yEmu) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = Get-ShapeById $slide $shapeId
    $shape.Left = ($xEmu / 12700.0) + $HalfEmuInPoints
    $shape.Top = ($yEmu / 12700.0) + $HalfEmuInPoints
}

# Slide 2
Move-Shape 2 44 7340303 849073
Move-Shape 2 46 7419165 2083188
Move-Shape 2 3  4268970 4542466

# Slide 3
Move-Shape 3 5  7490933 4519292

# Slide 4
Move-Shape 4 20 8047129 3258846
Move-Shape 4 24 7571517 4292600
Move-Shape 4 44 7353300 914842
Move-Shape 4 3  3905360 4591709
Move-Shape 4 7  7366743 3578277

# Slide 9
Move-Shape 9 6  5004984 5618614
